# Ajustes etapa 1, desarrollo etapa 2 y 3.
# Se adjuntan notas en el archivo excel (columnas D:F de Hoja1 con el
# seguimiento / revision de cada pagina) y se formatea la fila de
# encabezado (D1:F1) a juego con el resto de la hoja.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Encabezados nuevos en la fila 1 (D1:F1), mismo relleno que las
#     filas "ETAPA" + negrita ---------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("C1:F1").HorizontalAlignment = 1

$ws.Range("D1").Value = "Revision Coctel"
$ws.Range("E1").Value = "Revision Carvajal"
$ws.Range("F1").Value = "Ajustes Realizados"

# --- Notas / seguimiento por pagina (columna D, E y F) --------------
$ws.Range("D2").Value = "OK"
$ws.Range("E2").Value = "OK"
$ws.Range("F2").Value = "OK"

$ws.Range("D3").Value = "OK"

$ws.Range("D4").Value = "OK"
$ws.Range("E4").Value = "OK"
$ws.Range("F4").Value = "No hay"

$ws.Range("D5").Value = "OK"

$ws.Range("D6").Value = "Esperando ajuste de orden en elementos"

$ws.Range("D7").Value = "OK"
$ws.Range("E7").Value = "OK"

$ws.Range("D8").Value = "Esperando ajuste de orden en elementos"

# Fila 9 = encabezado "ETAPA 2": se resalta toda la fila A9:F9
$ws.Range("A9").Copy()
$ws.Range("B9:F9").PasteSpecial(-4122)

$ws.Range("D10").Value = "OK"
$ws.Range("D11").Value = "OK - Sky Scraper OK -  Falta texto por confirmar estructura"
$ws.Range("D12").Value = "No se tiene acceso, favor poner boton"
$ws.Range("D13").Value = "OK pendiente ajuste estructura DIVS tanto del buscador como del contenido"
$ws.Range("D14").Value = "OK pendiente resvisar un poco estructura de textos"
$ws.Range("D15").Value = "OK"
$ws.Range("D16").Value = "OK pendiente estructura de textos"

# Fila 17 = encabezado "ETAPA 3": se resalta toda la fila A17:F17
$ws.Range("A17").Copy()
$ws.Range("B17:F17").PasteSpecial(-4122)

$ws.Range("D18").Value = "OK"
$ws.Range("D19").Value = "Revisar texto que esta suelto al interior de los divs"
$ws.Range("D20").Value = "OK funciona igual que Home área temática"
$ws.Range("D21").Value = "OK funciona con el general.css ya generado"
$ws.Range("D22").Value = "OK ajustar botones por primary y colocar los a en los teléfonos"
$ws.Range("D23").Value = "OK se ajusta con el CSS de Administración guía de proveedores - home."
$ws.Range("D24").Value = "OK pendiente ajuste textos superior"

# Fila 25 = encabezado "ETAPA 4": se resalta toda la fila A25:F25
$ws.Range("A25").Copy()
$ws.Range("B25:F25").PasteSpecial(-4122)

# Fila 31 = encabezado "ETAPA 5": se resalta toda la fila A31:F31
$ws.Range("A31").Copy()
$ws.Range("B31:F31").PasteSpecial(-4122)

# --- Ancho de las columnas nuevas ------------------------------------
$ws.Columns.Item(4).ColumnWidth = 59.43
$ws.Columns.Item(5).ColumnWidth = 13.93
$ws.Columns.Item(6).ColumnWidth = 14.75

# --- Vista de la hoja -------------------------------------------------
$ws.Range("E7").Select()
$excel.ActiveWindow.ScrollColumn = 2

$app.CutCopyMode = $false
